$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 46; existing rows 46-94 shift down to 47-95.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new weekly record.
$ws.Range("A46").Value = 3
$ws.Range("B46").Value = "Femacal de La Calera"
$ws.Range("C46").Value = "Coquimbo"
$ws.Range("D46").Value = 44539
$ws.Range("E46").Value = 5
$ws.Range("F46").Value = 100112052
$ws.Range("G46").Value = "Albahaca"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 130
$ws.Range("K46").Value = 4500
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = 4769
$ws.Range("N46").Value = '$/docena de matas'
$ws.Range("O46").Value = "Provincia de Quillota"
$ws.Range("P46").Value = 795
$ws.Range("Q46").Value = 6
$ws.Range("R46").Value = "Hortaliza"
